$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44211
$ws.Range("H2").Value = 'Cultivar XV región'
$ws.Range("I2").Value = 'Segunda'
$ws.Range("K2").Value = 4500
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = 4750
$ws.Range("N2").Value = '$/caja 10 kilos'
$ws.Range("O2").Value = 'Región de Arica y Parinacota'
$ws.Range("P2").Value = 475
$ws.Range("Q2").Value = 10
$ws.Range("D3").Value = 44433
$ws.Range("I3").Value = 'Segunda'
$ws.Range("D4").Value = 44433
$ws.Range("I4").Value = 'Tercera'
$ws.Range("J4").Value = 120
$ws.Range("M4").Value = 14500
$ws.Range("P4").Value = 806
$ws.Range("D5").Value = 44398
$ws.Range("H5").Value = 'Cultivar IV Región'
$ws.Range("K5").Value = 17000
$ws.Range("L5").Value = 18000
$ws.Range("M5").Value = 17500
$ws.Range("N5").Value = '$/bandeja 18 kilos'
$ws.Range("O5").Value = 'Provincia de Limarí'
$ws.Range("P5").Value = 972
$ws.Range("Q5").Value = 18
$ws.Range("D6").Value = 44398
$ws.Range("H6").Value = 'Cultivar IV Región'
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 16000
$ws.Range("M6").Value = 15500
$ws.Range("N6").Value = '$/bandeja 18 kilos'
$ws.Range("O6").Value = 'Provincia de Limarí'
$ws.Range("P6").Value = 861
$ws.Range("Q6").Value = 18
$ws.Range("D7").Value = 44412
$ws.Range("H7").Value = 'Cultivar IV Región'
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 150
$ws.Range("K7").Value = 17000
$ws.Range("L7").Value = 18000
$ws.Range("M7").Value = 17500
$ws.Range("N7").Value = '$/bandeja 18 kilos'
$ws.Range("O7").Value = 'Provincia de Limarí'
$ws.Range("P7").Value = 972
$ws.Range("Q7").Value = 18
$ws.Range("D8").Value = 44742
$ws.Range("J8").Value = 250
$ws.Range("K8").Value = 15000
$ws.Range("L8").Value = 16000
$ws.Range("M8").Value = 15500
$ws.Range("P8").Value = 861
$ws.Range("D9").Value = 44762
$ws.Range("H9").Value = 'Cultivar IV Región'
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value = 160
$ws.Range("K9").Value = 15000
$ws.Range("L9").Value = 16000
$ws.Range("M9").Value = 15500
$ws.Range("N9").Value = '$/bandeja 18 kilos'
$ws.Range("O9").Value = 'Provincia de Limarí'
$ws.Range("P9").Value = 861
$ws.Range("Q9").Value = 18
$ws.Range("D10").Value = 44526
$ws.Range("H10").Value = 'Cultivar XV región'
$ws.Range("I10").Value = 'Primera'
$ws.Range("K10").Value = 5000
$ws.Range("L10").Value = 5500
$ws.Range("M10").Value = 5250
$ws.Range("N10").Value = '$/caja 10 kilos'
$ws.Range("O10").Value = 'Región de Arica y Parinacota'
$ws.Range("P10").Value = 525
$ws.Range("Q10").Value = 10
$ws.Range("D11").Value = 44526
$ws.Range("H11").Value = 'Cultivar XV región'
$ws.Range("I11").Value = 'Segunda'
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 4000
$ws.Range("L11").Value = 4500
$ws.Range("M11").Value = 4250
$ws.Range("N11").Value = '$/caja 10 kilos'
$ws.Range("O11").Value = 'Región de Arica y Parinacota'
$ws.Range("P11").Value = 425
$ws.Range("Q11").Value = 10
$ws.Range("D12").Value = 44526
$ws.Range("H12").Value = 'Cultivar XV región'
$ws.Range("I12").Value = 'Tercera'
$ws.Range("J12").Value = 120
$ws.Range("K12").Value = 3000
$ws.Range("L12").Value = 3500
$ws.Range("M12").Value = 3250
$ws.Range("N12").Value = '$/caja 10 kilos'
$ws.Range("O12").Value = 'Región de Arica y Parinacota'
$ws.Range("P12").Value = 325
$ws.Range("Q12").Value = 10
$ws.Range("D13").Value = 44769
$ws.Range("J13").Value = 140
$ws.Range("D14").Value = 44748
$ws.Range("H14").Value = 'Cultivar IV Región'
$ws.Range("J14").Value = 250
$ws.Range("K14").Value = 17000
$ws.Range("L14").Value = 18000
$ws.Range("M14").Value = 17500
$ws.Range("N14").Value = '$/bandeja 18 kilos'
$ws.Range("O14").Value = 'Provincia de Limarí'
$ws.Range("P14").Value = 972
$ws.Range("Q14").Value = 18
$ws.Range("D15").Value = 44783
$ws.Range("H15").Value = 'Cultivar IV Región'
$ws.Range("J15").Value = 150
$ws.Range("K15").Value = 17000
$ws.Range("L15").Value = 18000
$ws.Range("M15").Value = 17500
$ws.Range("N15").Value = '$/bandeja 18 kilos'
$ws.Range("O15").Value = 'Provincia de Limarí'
$ws.Range("P15").Value = 972
$ws.Range("Q15").Value = 18
$ws.Range("D16").Value = 44755
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value = 160
$ws.Range("D17").Value = 44554
$ws.Range("H17").Value = 'Cultivar XV región'
$ws.Range("I17").Value = 'Primera'
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 5000
$ws.Range("L17").Value = 6000
$ws.Range("M17").Value = 5500
$ws.Range("N17").Value = '$/caja 10 kilos'
$ws.Range("O17").Value = 'Región de Arica y Parinacota'
$ws.Range("P17").Value = 550
$ws.Range("Q17").Value = 10
$ws.Range("D18").Value = 44377
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 17000
$ws.Range("L18").Value = 18000
$ws.Range("M18").Value = 17600
$ws.Range("P18").Value = 978
$ws.Range("D19").Value = 44757
$ws.Range("H19").Value = 'Cultivar XV región'
$ws.Range("J19").Value = 150
$ws.Range("K19").Value = 6000
$ws.Range("L19").Value = 6500
$ws.Range("M19").Value = 6250
$ws.Range("N19").Value = '$/caja 10 kilos'
$ws.Range("O19").Value = 'Región de Arica y Parinacota'
$ws.Range("P19").Value = 625
$ws.Range("Q19").Value = 10
$ws.Range("D20").Value = 44363
$ws.Range("H20").Value = 'Cultivar IV Región'
$ws.Range("J20").Value = 140
$ws.Range("K20").Value = 14000
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = 14500
$ws.Range("N20").Value = '$/bandeja 18 kilos'
$ws.Range("O20").Value = 'Provincia de Limarí'
$ws.Range("P20").Value = 806
$ws.Range("Q20").Value = 18
$ws.Range("D21").Value = 44776
$ws.Range("H21").Value = 'Cultivar IV Región'
$ws.Range("I21").Value = 'Primera'
$ws.Range("J21").Value = 200
$ws.Range("K21").Value = 17000
$ws.Range("L21").Value = 18000
$ws.Range("M21").Value = 17500
$ws.Range("N21").Value = '$/bandeja 18 kilos'
$ws.Range("O21").Value = 'Provincia de Limarí'
$ws.Range("P21").Value = 972
$ws.Range("Q21").Value = 18
$ws.Range("D22").Value = 44391
$ws.Range("H22").Value = 'Cultivar IV Región'
$ws.Range("I22").Value = 'Segunda'
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 15000
$ws.Range("L22").Value = 16000
$ws.Range("M22").Value = 15500
$ws.Range("N22").Value = '$/bandeja 18 kilos'
$ws.Range("O22").Value = 'Provincia de Limarí'
$ws.Range("P22").Value = 861
$ws.Range("Q22").Value = 18
$ws.Range("D23").Value = 44533
$ws.Range("H23").Value = 'Cultivar XV región'
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = 6000
$ws.Range("L23").Value = 7000
$ws.Range("M23").Value = 6500
$ws.Range("N23").Value = '$/caja 10 kilos'
$ws.Range("O23").Value = 'Región de Arica y Parinacota'
$ws.Range("P23").Value = 650
$ws.Range("Q23").Value = 10
$ws.Range("D24").Value = 44533
$ws.Range("H24").Value = 'Cultivar XV región'
$ws.Range("I24").Value = 'Segunda'
$ws.Range("J24").Value = 120
$ws.Range("K24").Value = 4000
$ws.Range("L24").Value = 5000
$ws.Range("M24").Value = 4500
$ws.Range("N24").Value = '$/caja 10 kilos'
$ws.Range("O24").Value = 'Región de Arica y Parinacota'
$ws.Range("P24").Value = 450
$ws.Range("Q24").Value = 10
$ws.Range("D25").Value = 44454
$ws.Range("J25").Value = 160
$ws.Range("K25").Value = 19000
$ws.Range("L25").Value = 20000
$ws.Range("M25").Value = 19500
$ws.Range("P25").Value = 1083
$ws.Range("D26").Value = 44405
$ws.Range("I26").Value = 'Segunda'
$ws.Range("D27").Value = 44771
$ws.Range("H27").Value = 'Cultivar XV región'
$ws.Range("J27").Value = 140
$ws.Range("K27").Value = 8000
$ws.Range("L27").Value = 9000
$ws.Range("M27").Value = 8500
$ws.Range("N27").Value = '$/caja 10 kilos'
$ws.Range("O27").Value = 'Región de Arica y Parinacota'
$ws.Range("P27").Value = 850
$ws.Range("Q27").Value = 10
$ws.Range("D28").Value = 44221
$ws.Range("H28").Value = 'Cultivar XV región'
$ws.Range("J28").Value = 140
$ws.Range("K28").Value = 5000
$ws.Range("L28").Value = 6000
$ws.Range("M28").Value = 5500
$ws.Range("N28").Value = '$/caja 10 kilos'
$ws.Range("O28").Value = 'Región de Arica y Parinacota'
$ws.Range("P28").Value = 550
$ws.Range("Q28").Value = 10
$ws.Range("D29").Value = 44435
$ws.Range("H29").Value = 'Cultivar IV Región'
$ws.Range("I29").Value = 'Segunda'
$ws.Range("J29").Value = 100
$ws.Range("K29").Value = 17000
$ws.Range("L29").Value = 18000
$ws.Range("M29").Value = 17500
$ws.Range("N29").Value = '$/bandeja 18 kilos'
$ws.Range("O29").Value = 'Provincia de Limarí'
$ws.Range("P29").Value = 972
$ws.Range("Q29").Value = 18
$ws.Range("D30").Value = 44435
$ws.Range("I30").Value = 'Tercera'
$ws.Range("J30").Value = 120
$ws.Range("M30").Value = 14500
$ws.Range("P30").Value = 806
